$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell 2 4 "42.597.94"
Set-TextCell 2 5 "  +1.30%  "

# Row 3
Set-TextCell 3 4 "2.251.44"
Set-TextCell 3 5 "  +0.51%  "

# Row 5
Set-TextCell 5 4 "246.62"
Set-TextCell 5 5 "  -0.06%  "

# Row 6
Set-TextCell 6 4 "0.631"
Set-TextCell 6 5 "  -0.05%  "

# Row 7
Set-TextCell 7 4 "76.43"
Set-TextCell 7 5 "  +0.74%  "

# Row 8
Set-TextCell 8 5 "  +0.15%  "

# Row 9
Set-TextCell 9 4 "0.627"
Set-TextCell 9 5 "  -0.24%  "

# Row 10
Set-TextCell 10 4 "44.61"
Set-TextCell 10 5 "  +10.18%  "

# Row 11
Set-TextCell 11 4 "0.0954"
Set-TextCell 11 5 "  +0.01%  "

# Row 12
Set-TextCell 12 4 "7.35"
Set-TextCell 12 5 "  +3.18%  "

# Row 13
Set-TextCell 13 5 "  -0.42%  "

# Row 14
Set-TextCell 14 4 "14.75"
Set-TextCell 14 5 "  -0.76%  "

# Row 15
Set-TextCell 15 4 "0.866"
Set-TextCell 15 5 "  +0.43%  "

# Row 16
Set-TextCell 16 4 "2.266.30"
Set-TextCell 16 5 "  +0.50%  "

# Row 17
Set-TextCell 17 4 "42.420.19"
Set-TextCell 17 5 "  +1.19%  "

# Row 18
Set-TextCell 18 5 "  +3.62%  "

# Row 19
Set-TextCell 19 4 "6.22"
Set-TextCell 19 5 "  +1.26%  "

# Row 20
Set-TextCell 20 4 "72.26"
Set-TextCell 20 5 "  +0.85%  "

# Row 21
Set-TextCell 21 4 "10.96"
Set-TextCell 21 5 "  +53.00%  "

# Row 22
Set-TextCell 22 5 "  -0.07%  "

# Row 23
Set-TextCell 23 4 "232.48"
Set-TextCell 23 5 "  +1.08%  "

# Row 24
Set-TextCell 24 4 "11.82"
Set-TextCell 24 5 "  +3.11%  "

# Row 25
Set-TextCell 25 5 "  +0.06%  "

# Row 26
Set-TextCell 26 5 "  -1.56%  "

# Row 27
Set-TextCell 27 5 "  +0.16%  "

# Row 28
Set-TextCell 28 5 "  +3.96%  "

# Row 29
Set-TextCell 29 4 "167.42"
Set-TextCell 29 5 "  -0.85%  "

# Row 30
Set-TextCell 30 4 "20.73"
Set-TextCell 30 5 "  +0.89%  "

# Row 31
Set-TextCell 31 2 "Filecoin"
Set-TextCell 31 3 "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell 31 4 "5.77"
Set-TextCell 31 5 "  +18.89%  "

# Row 32
Set-TextCell 32 4 "0.0824"
Set-TextCell 32 5 "  -2.48%  "

# Row 33
Set-TextCell 33 2 "InjectiveProtocol"
Set-TextCell 33 3 "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell 33 4 "32.14"
Set-TextCell 33 5 "  -3.93%  "

# Row 34
Set-TextCell 34 4 "0.120"
Set-TextCell 34 5 "  -0.93%  "

# Row 35
Set-TextCell 35 2 "RenderToken"
Set-TextCell 35 3 "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell 35 4 "4.81"
Set-TextCell 35 5 "  +6.47%  "

# Row 36
Set-TextCell 36 2 "Stellar"
Set-TextCell 36 3 "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell 36 4 "0.126"
Set-TextCell 36 5 "  -0.02%  "

# Row 37
Set-TextCell 37 4 "0.0317"
Set-TextCell 37 5 "  +5.79%  "

# Row 38
Set-TextCell 38 4 "14.25"
Set-TextCell 38 5 "  +6.36%  "

# Row 39
Set-TextCell 39 5 "  +0.76%  "

# Row 40
Set-TextCell 40 4 "5.81"
Set-TextCell 40 5 "  -2.11%  "

# Row 41
Set-TextCell 41 4 "64.42"
Set-TextCell 41 5 "  +6.35%  "

# Row 42
Set-TextCell 42 5 "  +0.12%  "

# Row 43
Set-TextCell 43 4 "108.30"
Set-TextCell 43 5 "  -2.89%  "

# Row 44
Set-TextCell 44 4 "8.93"
Set-TextCell 44 5 "  +2.31%  "

# Row 45
Set-TextCell 45 5 "  +2.13%  "

# Row 46
Set-TextCell 46 4 "0.997"
Set-TextCell 46 5 "  +0.08%  "

# Row 47
Set-TextCell 47 2 "NEARProtocol"
Set-TextCell 47 3 "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell 47 4 "2.41"
Set-TextCell 47 5 "  +6.99%  "

# Row 48
Set-TextCell 48 2 "ARBITRUM"
Set-TextCell 48 3 "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell 48 4 "1.15"
Set-TextCell 48 5 "  +0.80%  "

# Row 49
Set-TextCell 49 4 "1.19"
Set-TextCell 49 5 "  +1.87%  "

# Row 50
Set-TextCell 50 4 "4.15"

# Row 51
Set-TextCell 51 2 "HuobiToken"
Set-TextCell 51 3 "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextCell 51 4 "2.71"
Set-TextCell 51 5 "  +0.90%  "
